$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test_number column (A) from 2 to 3 for the data rows (2-11)
$ws.Range("A2:A11").Value = 3

# Move the active selection to P7 (as last saved by the author)
$ws.Range("P7").Select()
